# Listas sem duplicação de professores
# Replace cells that contained a list of teacher-occupied slots (with
# duplicated teacher names across multiple classes) with a plain "-"
# to avoid showing the same teacher listed more than once.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @("D2", "E2", "D8", "E8", "F11", "F12", "F14", "F15", "E18", "F18", "E19", "F19", "C21", "E21")

foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).Value = "-"
}
